# Insert a new price-report row at row 615 ("Vega Modelo de Temuco" /
# "Zapallo italiano" weekly update), pushing the existing rows 615-678
# down to 616-679.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 615..678 down by one (inherits formatting from the row above,
# matching the existing "s=2" date style already used throughout column D).
$ws.Rows.Item(615).Insert()

# Populate the newly inserted row 615 with the new weekly record.
$ws.Cells.Item(615, 1).Value = 10
$ws.Cells.Item(615, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(615, 3).Value = "La Araucanía"
$ws.Cells.Item(615, 4).Value = 44946
$ws.Cells.Item(615, 5).Value = 9
$ws.Cells.Item(615, 6).Value = 100112032
$ws.Cells.Item(615, 7).Value = "Zapallo italiano"
$ws.Cells.Item(615, 8).Value = "Sin especificar"
$ws.Cells.Item(615, 9).Value = "Primera"
$ws.Cells.Item(615, 10).Value = 155
$ws.Cells.Item(615, 11).Value = 10000
$ws.Cells.Item(615, 12).Value = 10000
$ws.Cells.Item(615, 13).Value = 10000
$ws.Cells.Item(615, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(615, 15).Value = "Región del Maule"
$ws.Cells.Item(615, 16).Value = 200
$ws.Cells.Item(615, 17).Value = 50
$ws.Cells.Item(615, 18).Value = "Hortaliza"
